$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "10.00", "42.612.48") are preserved exactly as text, matching the
# source data which stores these as inline/shared strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.612.48"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.527.45"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "315.09"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "98.99"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "35.22"
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "7.21"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").Value = "2.911.71"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "15.25"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "2.516.70"
$ws.Range("E16").Value = "  -3.16%  "
$ws.Range("D17").Value = "0.810"
$ws.Range("E17").Value = "  -4.41%  "
$ws.Range("D18").Value = "42.585.60"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "6.60"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "12.20"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0$([char]0x2083)0939"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "69.07"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "242.01"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").Value = "2.86"
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "25.58"
$ws.Range("E27").Value = "  -3.96%  "
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("D29").Value = "10.00"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "37.73"
$ws.Range("E30").Value = "  -8.26%  "
$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "156.30"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").Value = "17.53"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "22.03"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0295"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.016.61"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").Value = "3.22"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "8.96"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "2.763.24"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "79.19"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").Value = "71.48"
$ws.Range("E51").Value = "  -3.12%  "
